$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "value" header (B1) to "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# Propagate the date style (from A2) down through the new rows (A3:A22)
$ws.Range("A2").Copy($ws.Range("A3:A22"))

# Clear the old B2 value since row 2 now only holds a date (no value)
$ws.Range("B2").ClearContents()

# Fill in the date / value series (rows shifted down by one, with a
# new first date row and a new last date-only row)
$ws.Range("A2").Value = 38717

$ws.Range("A3").Value = 39082
$ws.Range("B3").Value = 6.681483765882756

$ws.Range("A4").Value = 39447
$ws.Range("B4").Value = 5.732148352530309

$ws.Range("A5").Value = 39813
$ws.Range("B5").Value = 6.181322443148352

$ws.Range("A6").Value = 40178
$ws.Range("B6").Value = 9.97031398925483

$ws.Range("A7").Value = 40543
$ws.Range("B7").Value = 5.968279190641868

$ws.Range("A8").Value = 40908
$ws.Range("B8").Value = 6.387913216057295

$ws.Range("A9").Value = 41274
$ws.Range("B9").Value = 3.489647115587391

$ws.Range("A10").Value = 41639
$ws.Range("B10").Value = 3.062667370145955

$ws.Range("A11").Value = 42004
$ws.Range("B11").Value = 1.40861416720266

$ws.Range("A12").Value = 42369
$ws.Range("B12").Value = 1.975538030067248

$ws.Range("A13").Value = 42735
$ws.Range("B13").Value = 2.149250550875026

$ws.Range("A14").Value = 43100
$ws.Range("B14").Value = 2.516312190944614

$ws.Range("A15").Value = 43465
$ws.Range("B15").Value = 2.494967260739056

$ws.Range("A16").Value = 43830
$ws.Range("B16").Value = 0.9136132777513017

$ws.Range("A17").Value = 44196
$ws.Range("B17").Value = 2.275661779503824

$ws.Range("A18").Value = 44561
$ws.Range("B18").Value = 3.546865287857126

$ws.Range("A19").Value = 44926
$ws.Range("B19").Value = 4.646251873334628

$ws.Range("A20").Value = 45291
$ws.Range("B20").Value = 2.244754177395403

$ws.Range("A21").Value = 45657
$ws.Range("B21").Value = 1.608247521160311

$ws.Range("A22").Value = 46022
